$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 16 data, mirroring the style/content pattern of row 15
$ws.Range("A16").Value = 14

# Copy A15's formatting (bold, centered, bordered) onto A16 without
# disturbing its value
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)

$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 0.9872764148807824
$ws.Range("D16").Value = 1.033312647539365
$ws.Range("E16").Value = 0.9838123690021011
$ws.Range("F16").Value = 0.9872764148807824
$ws.Range("G16").Value = 1.015657891960751
$ws.Range("H16").Value = 0.9651918068774645
$ws.Range("I16").Value = 0.9848434837540447
$ws.Range("J16").Value = 1.033312647539365
$ws.Range("K16").Value = 1.008562508270733
$ws.Range("L16").Value = 0.9979194615757578
$ws.Range("M16").Value = 0.9950157690024182
